$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("version")

# Row 22: nghttp2 section -> bump msvc15 build tag v1.33.0 -> v1.34.0 and add changelog line
$ws.Range("B22").Value = "msvc15 / msvc15-v1.34.0"
$ws.Range("C22").Value = "v1.34"
$ws.Range("D22").Value = "0206 master => v1.30.0`n0301 v1.30.0 => v1.31.0`n0415 v1.31.0 => v1.31.1`n0514 v1.32.0`n1006 v1.34.0`n0911 v1.33.0"

# Row 39: subversion section -> bump 1.10.2 -> 1.10.3 and add changelog line
$ws.Range("C39").Value = "1.10.3"
$ws.Range("D39").Value = "0328 1.10.0-rc1`n0410 1.10.0-rc1 => 1.10.0-rc2`n0411 1.10.0-rc2 => 1.10.0`n0724 1.10.2`n1006 1.10.3"

# Update row heights to reflect the extra wrapped line (autofit)
$ws.Rows.Item(22).RowHeight = 90
$ws.Rows.Item(39).RowHeight = 75

# Update the selection/view state to match final editing position
$ws.Range("B28").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 28
